$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header swap: BP1/BQ1 (average_doctor <-> average_doctor_old) ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Updated statistic values (harvard case classification) ---
$ws.Range("AI4").Value = 0.224
$ws.Range("AU4").Value = 0.15
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.168
$ws.Range("BA4").Value = 1.987
$ws.Range("BB4").Value = 0.158
$ws.Range("BC4").Value = 0.397
$ws.Range("BG4").Value = 0.739
$ws.Range("BH4").Value = 0.131
$ws.Range("BI4").Value = 0.361
$ws.Range("BM4").Value = 0.706
$ws.Range("BN4").Value = 0.086
$ws.Range("BO4").Value = 0.293
$ws.Range("BP4").Value = 0.662
$ws.Range("BQ4").Value = 0.663
$ws.Range("E4").Value = 0.404
$ws.Range("F4").Value = 0.085
$ws.Range("G4").Value = 0.291
$ws.Range("N4").Value = 0.421
$ws.Range("O4").Value = 0.068
$ws.Range("P4").Value = 0.26
$ws.Range("W4").Value = 0.231
$ws.Range("X4").Value = 0.105
$ws.Range("Y4").Value = 0.324
$ws.Range("AI5").Value = 0.267
$ws.Range("AU5").Value = 0.304
$ws.Range("AV5").Value = 0.104
$ws.Range("AW5").Value = 0.322
$ws.Range("BA5").Value = 1.371
$ws.Range("BG5").Value = 0.409
$ws.Range("BH5").Value = 0.043
$ws.Range("BI5").Value = 0.208
$ws.Range("BM5").Value = 0.561
$ws.Range("BN5").Value = 0.068
$ws.Range("BO5").Value = 0.261
$ws.Range("BP5").Value = 0.457
$ws.Range("BQ5").Value = 0.451
$ws.Range("E5").Value = 0.5
$ws.Range("F5").Value = 0.09
$ws.Range("G5").Value = 0.3
$ws.Range("N5").Value = 0.748
$ws.Range("O5").Value = 0.079
$ws.Range("P5").Value = 0.281
$ws.Range("W5").Value = 0.233
$ws.Range("X5").Value = 0.109
$ws.Range("Y5").Value = 0.331
$ws.Range("AI6").Value = 0.244
$ws.Range("AU6").Value = 0.201
$ws.Range("BA6").Value = 1.613
$ws.Range("BG6").Value = 0.527
$ws.Range("BM6").Value = 0.625
$ws.Range("BP6").Value = 0.538
$ws.Range("BQ6").Value = 0.534
$ws.Range("E6").Value = 0.447
$ws.Range("N6").Value = 0.539
$ws.Range("W6").Value = 0.232
$ws.Range("AI7").Value = 0.257
$ws.Range("AU7").Value = 0.252
$ws.Range("BA7").Value = 1.457
$ws.Range("BG7").Value = 0.449
$ws.Range("BM7").Value = 0.585
$ws.Range("BP7").Value = 0.486
$ws.Range("BQ7").Value = 0.48
$ws.Range("E7").Value = 0.477
$ws.Range("N7").Value = 0.647
$ws.Range("W7").Value = 0.233
$ws.Range("AI8").Value = 0.245
$ws.Range("AJ8").Value = 0.097
$ws.Range("AK8").Value = 0.312
$ws.Range("AU8").Value = 0.241
$ws.Range("AV8").Value = 0.073
$ws.Range("AW8").Value = 0.271
$ws.Range("BA8").Value = 1.71
$ws.Range("BB8").Value = 0.13
$ws.Range("BC8").Value = 0.36
$ws.Range("BG8").Value = 0.563
$ws.Range("BH8").Value = 0.096
$ws.Range("BI8").Value = 0.309
$ws.Range("BM8").Value = 0.693
$ws.Range("BN8").Value = 0.064
$ws.Range("BO8").Value = 0.253
$ws.Range("BP8").Value = 0.57
$ws.Range("BQ8").Value = 0.579
$ws.Range("E8").Value = 0.537
$ws.Range("F8").Value = 0.118
$ws.Range("G8").Value = 0.343
$ws.Range("N8").Value = 0.753
$ws.Range("O8").Value = 0.066
$ws.Range("P8").Value = 0.257
$ws.Range("W8").Value = 0.234
$ws.Range("X8").Value = 0.111
$ws.Range("Y8").Value = 0.333
$ws.Range("AI9").Value = 0.136
$ws.Range("AJ9").Value = 0.118
$ws.Range("AK9").Value = 0.343
$ws.Range("BA9").Value = 1.614
$ws.Range("BB9").Value = 0.242
$ws.Range("BC9").Value = 0.492
$ws.Range("BG9").Value = 0.591
$ws.Range("BH9").Value = 0.242
$ws.Range("BI9").Value = 0.492
$ws.Range("BM9").Value = 0.614
$ws.Range("BN9").Value = 0.237
$ws.Range("BO9").Value = 0.487
$ws.Range("BP9").Value = 0.538
$ws.Range("BQ9").Value = 0.534
$ws.Range("E9").Value = 0.455
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.636
$ws.Range("O9").Value = 0.231
$ws.Range("P9").Value = 0.481
$ws.Range("W9").Value = 0.136
$ws.Range("X9").Value = 0.118
$ws.Range("Y9").Value = 0.343
$ws.Range("AI10").Value = 0.273
$ws.Range("AJ10").Value = 0.198
$ws.Range("AK10").Value = 0.445
$ws.Range("AU10").Value = 0.227
$ws.Range("AV10").Value = 0.176
$ws.Range("AW10").Value = 0.419
$ws.Range("BA10").Value = 1.977
$ws.Range("BB10").Value = 0.249
$ws.Range("BC10").Value = 0.499
$ws.Range("BG10").Value = 0.636
$ws.Range("BH10").Value = 0.231
$ws.Range("BI10").Value = 0.481
$ws.Range("BM10").Value = 0.864
$ws.Range("BN10").Value = 0.118
$ws.Range("BO10").Value = 0.343
$ws.Range("BP10").Value = 0.659
$ws.Range("BQ10").Value = 0.679
$ws.Range("E10").Value = 0.591
$ws.Range("F10").Value = 0.242
$ws.Range("G10").Value = 0.492
$ws.Range("N10").Value = 0.841
$ws.Range("O10").Value = 0.134
$ws.Range("P10").Value = 0.366
$ws.Range("W10").Value = 0.273
$ws.Range("X10").Value = 0.198
$ws.Range("Y10").Value = 0.445
$ws.Range("AI11").Value = 0.273
$ws.Range("AJ11").Value = 0.198
$ws.Range("AK11").Value = 0.445
$ws.Range("AU11").Value = 0.341
$ws.Range("AV11").Value = 0.225
$ws.Range("AW11").Value = 0.474
$ws.Range("BA11").Value = 1.977
$ws.Range("BB11").Value = 0.249
$ws.Range("BC11").Value = 0.499
$ws.Range("BG11").Value = 0.636
$ws.Range("BH11").Value = 0.231
$ws.Range("BI11").Value = 0.481
$ws.Range("BM11").Value = 0.864
$ws.Range("BN11").Value = 0.118
$ws.Range("BO11").Value = 0.343
$ws.Range("BP11").Value = 0.659
$ws.Range("BQ11").Value = 0.679
$ws.Range("E11").Value = 0.614
$ws.Range("F11").Value = 0.237
$ws.Range("G11").Value = 0.487
$ws.Range("N11").Value = 0.864
$ws.Range("O11").Value = 0.118
$ws.Range("P11").Value = 0.343
$ws.Range("W11").Value = 0.273
$ws.Range("X11").Value = 0.198
$ws.Range("Y11").Value = 0.445
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.667
$ws.Range("AV12").Value = 1.689
$ws.Range("AW12").Value = 1.3
$ws.Range("BA12").Value = 3.704
$ws.Range("BB12").Value = 0.372
$ws.Range("BC12").Value = 0.61
$ws.Range("BG12").Value = 1.071
$ws.Range("BH12").Value = 0.066
$ws.Range("BI12").Value = 0.258
$ws.Range("BM12").Value = 1.395
$ws.Range("BN12").Value = 0.449
$ws.Range("BO12").Value = 0.67
$ws.Range("BP12").Value = 1.235
$ws.Range("BQ12").Value = 1.28
$ws.Range("E12").Value = 1.444
$ws.Range("F12").Value = 0.691
$ws.Range("G12").Value = 0.831
$ws.Range("N12").Value = 1.625
$ws.Range("O12").Value = 1.534
$ws.Range("P12").Value = 1.239
$ws.Range("W12").Value = 1.75
$ws.Range("X12").Value = 0.688
$ws.Range("Y12").Value = 0.829
$ws.Range("AI13").Value = 1.365
$ws.Range("AJ13").Value = 0.41
$ws.Range("AK13").Value = 0.64
$ws.Range("AU13").Value = 2.35
$ws.Range("AV13").Value = 0.749
$ws.Range("AW13").Value = 0.865
$ws.Range("BA13").Value = 2.48
$ws.Range("BB13").Value = 0.313
$ws.Range("BC13").Value = 0.56
$ws.Range("BG13").Value = 0.619
$ws.Range("BH13").Value = 0.089
$ws.Range("BI13").Value = 0.298
$ws.Range("BM13").Value = 0.953
$ws.Range("BN13").Value = 0.376
$ws.Range("BO13").Value = 0.613
$ws.Range("BP13").Value = 0.827
$ws.Range("BQ13").Value = 0.762
$ws.Range("E13").Value = 1.642
$ws.Range("F13").Value = 0.853
$ws.Range("G13").Value = 0.923
$ws.Range("N13").Value = 2.202
$ws.Range("O13").Value = 0.96
$ws.Range("P13").Value = 0.98
$ws.Range("W13").Value = 1.076
$ws.Range("X13").Value = 0.176
$ws.Range("Y13").Value = 0.42
